$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formula in B2 with a static value (support for updates as well)
$ws.Range("B2").Value = "relay_100"

# Update the active selection to match the new state
$ws.Range("F15").Select()
